# Add a new slide 4 using the "Title and Content" layout (ppLayoutObject = 16),
# matching the other content slides in this deck (slideLayout2.xml).
$p = $ppt.ActivePresentation
$s = $p.Slides.Add(4, 16)

# Set the title text; leave the body/content placeholder empty, as in the target slide.
$s.Shapes.Title.TextFrame.TextRange.Text = "Adaptive Rule Results"
